# Append the two new trading-day rows (2025-11-24 / 2025-11-25, serials
# 45985 / 45986) to every "date"/"remn_amt" worksheet in the workbook.
# Each worksheet already ends with a contiguous block of date/value rows
# starting at row 2 (row 1 is the header); we just continue that block.

$wb = $excel.ActiveWorkbook

# New rows to add per worksheet (by 1-based worksheet index), matching the
# order of sheets in the workbook: 카카오, NAVER, 농심, 삼양식품, 엔씨소프트.
$newRows = @{
    1 = @(
        @{ Date = 45985; Value = 786542 },
        @{ Date = 45986; Value = 776488 }
    )
    2 = @(
        @{ Date = 45985; Value = 1322255 },
        @{ Date = 45986; Value = 1303364 }
    )
    3 = @(
        @{ Date = 45985; Value = 111509 },
        @{ Date = 45986; Value = 114145 }
    )
    4 = @(
        @{ Date = 45985; Value = 349793 },
        @{ Date = 45986; Value = 367451 }
    )
    5 = @(
        @{ Date = 45985; Value = 238474 },
        @{ Date = 45986; Value = 264652 }
    )
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Find the last already-populated row in column A (the date column).
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

    $rows = $newRows[$i]
    foreach ($row in $rows) {
        $lastRow = $lastRow + 1
        $ws.Cells.Item($lastRow, 1).Value = $row.Date
        $ws.Cells.Item($lastRow, 2).Value = $row.Value
        # Preserve the date-time number format used by the rest of column A.
        $ws.Cells.Item($lastRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
}
